$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4295.4863
$ws.Range("I11").Value = 4295.4863
$ws.Range("K11").Value = 4295.4863
$ws.Range("M11").Value = -4155.4863

$ws.Range("H17").Value = 4086.3635
$ws.Range("J17").Value = 4086.3635
$ws.Range("L17").Value = 12259.0905
$ws.Range("N17").Value = -12595.0905

$ws.Range("H82").Value = 934.3333
$ws.Range("I82").Value = 988.625
$ws.Range("J82").Value = 500
$ws.Range("K82").Value = 2965.875
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -2559.875
$ws.Range("N82").Value = -2312

$ws.Range("H85").Value = 934.3333
$ws.Range("I85").Value = 988.625
$ws.Range("J85").Value = 500
$ws.Range("K85").Value = 2965.875
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = -1561.875
$ws.Range("N85").Value = -4308

$ws.Range("H92").Value = 1475.2759
$ws.Range("I92").Value = 1076.091
$ws.Range("J92").Value = 2729.8572
$ws.Range("K92").Value = 1076.091
$ws.Range("L92").Value = 2729.8572
$ws.Range("M92").Value = 171.9090000000001
$ws.Range("N92").Value = -5225.8572

$ws.Range("H132").Value = 1190.7632
$ws.Range("I132").Value = 745.48486
$ws.Range("J132").Value = 4129.6
$ws.Range("K132").Value = 2236.45458
$ws.Range("L132").Value = 12388.8
$ws.Range("M132").Value = 293.5454199999999
$ws.Range("N132").Value = -17448.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 2000
$ws.Range("K14").Value = 2000
$ws.Range("M14").Value = -1825

$ws.Range("H74").Value = 5612.0835
$ws.Range("I74").Value = 3499.5
$ws.Range("J74").Value = 6034.6
$ws.Range("K74").Value = 3499.5
$ws.Range("L74").Value = 6034.6
$ws.Range("M74").Value = -2625.5
$ws.Range("N74").Value = -7782.6

$ws.Range("H77").Value = 5612.0835
$ws.Range("I77").Value = 3499.5
$ws.Range("J77").Value = 6034.6
$ws.Range("K77").Value = 17497.5
$ws.Range("L77").Value = 30173
$ws.Range("M77").Value = -13129.5
$ws.Range("N77").Value = -38909

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H107").Value = 1481.5
$ws.Range("I107").Value = 1467.2778
$ws.Range("K107").Value = 1467.2778
$ws.Range("M107").Value = 452.7221999999999

$ws.Range("H134").Value = 8873.875
$ws.Range("I134").Value = 7915.6924
$ws.Range("J134").Value = 10006.272
$ws.Range("K134").Value = 23747.0772
$ws.Range("L134").Value = 30018.816
$ws.Range("M134").Value = -21212.0772
$ws.Range("N134").Value = -35088.81600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7133.6665
$ws.Range("I86").Value = 5933.3335
$ws.Range("J86").Value = 8334
$ws.Range("K86").Value = 5933.3335
$ws.Range("L86").Value = 8334
$ws.Range("M86").Value = -4810.3335
$ws.Range("N86").Value = -10580

$ws.Range("H89").Value = 7133.6665
$ws.Range("I89").Value = 5933.3335
$ws.Range("J89").Value = 8334
$ws.Range("K89").Value = 29666.6675
$ws.Range("L89").Value = 41670
$ws.Range("M89").Value = -24050.6675
$ws.Range("N89").Value = -52902

$ws.Range("H94").Value = 2828.158
$ws.Range("I94").Value = 8181
$ws.Range("J94").Value = 916.4286
$ws.Range("K94").Value = 8181
$ws.Range("L94").Value = 916.4286
$ws.Range("M94").Value = -7730
$ws.Range("N94").Value = -1818.4286

$ws.Range("H134").Value = 12340.286
$ws.Range("I134").Value = 11972
$ws.Range("K134").Value = 35916
$ws.Range("M134").Value = -33381

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2423.1667
$ws.Range("I58").Value = 1517.5
$ws.Range("J58").Value = 2876
$ws.Range("K58").Value = 4552.5
$ws.Range("L58").Value = 8628
$ws.Range("M58").Value = -4424.5
$ws.Range("N58").Value = -8884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4847.125
$ws.Range("I80").Value = 3999.75
$ws.Range("K80").Value = 3999.75
$ws.Range("M80").Value = -3001.75

$ws.Range("H83").Value = 4847.125
$ws.Range("I83").Value = 3999.75
$ws.Range("K83").Value = 19998.75
$ws.Range("M83").Value = -15006.75

$ws.Range("H126").Value = 2639.8965
$ws.Range("I126").Value = 1920.8636
$ws.Range("K126").Value = 5762.5908
$ws.Range("M126").Value = -3292.5908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3499.5
$ws.Range("I7").Value = 3499.5
$ws.Range("K7").Value = 3499.5
$ws.Range("M7").Value = -3387.5

$ws.Range("H16").Value = 624.75
$ws.Range("I16").Value = 565.16
$ws.Range("J16").Value = 1121.3334
$ws.Range("K16").Value = 565.16
$ws.Range("L16").Value = 1121.3334
$ws.Range("M16").Value = -395.16
$ws.Range("N16").Value = -1461.3334

$ws.Range("H61").Value = 4095.0715
$ws.Range("I61").Value = 3740.875
$ws.Range("K61").Value = 3740.875
$ws.Range("M61").Value = -3538.875

$ws.Range("H113").Value = 4095.0715
$ws.Range("I113").Value = 3740.875
$ws.Range("K113").Value = 3740.875
$ws.Range("M113").Value = -1570.875

$ws.Range("H126").Value = 3499.5
$ws.Range("I126").Value = 3499.5
$ws.Range("K126").Value = 10498.5
$ws.Range("M126").Value = -8028.5

$ws.Range("H132").Value = 5203.2856
$ws.Range("I132").Value = 2624.8572
$ws.Range("J132").Value = 7781.7144
$ws.Range("K132").Value = 7874.571599999999
$ws.Range("L132").Value = 23345.1432
$ws.Range("M132").Value = -5344.571599999999
$ws.Range("N132").Value = -28405.1432

$ws.Range("H135").Value = 104436.2
$ws.Range("J135").Value = 104436.2
$ws.Range("L135").Value = 104436.2
$ws.Range("N135").Value = -114576.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3791.2058
$ws.Range("J81").Value = 6636.8
$ws.Range("L81").Value = 13273.6
$ws.Range("N81").Value = -15395.6

$ws.Range("H84").Value = 3791.2058
$ws.Range("J84").Value = 6636.8
$ws.Range("L84").Value = 66368
$ws.Range("N84").Value = -76976

$ws.Range("H96").Value = 1598.6666
$ws.Range("J96").Value = 1649
$ws.Range("L96").Value = 1649
$ws.Range("N96").Value = -4395

$ws.Range("H100").Value = 1573.1333
$ws.Range("I100").Value = 1335.2858
$ws.Range("K100").Value = 2670.5716
$ws.Range("M100").Value = -2129.5716

$ws.Range("H132").Value = 6607.7393
$ws.Range("I132").Value = 5012.067
$ws.Range("K132").Value = 15036.201
$ws.Range("M132").Value = -12506.201
